$wb = $excel.ActiveWorkbook

# --- "Modify Transaction" sheet (sheet4): insert a ReceiptNumber row before the
#     existing "submitmakerepayment" row, pushing it down from row 4 to row 5.
$wsMT1 = $wb.Worksheets.Item("Modify Transaction")
$wsMT1.Rows.Item(4).Insert()
$wsMT1.Range("A4").Value = "ReceiptNumber"
$wsMT1.Range("B4").Value = 12345

# --- "Modify Transaction1" sheet (sheet5): same kind of insert, different value.
$wsMT2 = $wb.Worksheets.Item("Modify Transaction1")
$wsMT2.Rows.Item(4).Insert()
$wsMT2.Range("A4").Value = "ReceiptNumber"
$wsMT2.Range("B4").Value = 7654

# --- Update the selections shown on each relevant sheet.
$wsMT1.Range("B8").Select()
$wsMT2.Range("B8").Select()

# --- "Transactions" sheet (sheet8) becomes the active/selected tab.
$wsTxn = $wb.Worksheets.Item("Transactions")
$wsTxn.Range("H8").Select()
